$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update result rows 2-25 (columns B:N, skipping G/K which remain 0)
# with newly computed values for the "case with 380 kV" scenario
# (pl_mw.xlsx results table).

$ws.Range("B2").Value = 2.767168770399508
$ws.Range("C2").Value = 0.3573878827432679
$ws.Range("D2").Value = 0.009189471803441762
$ws.Range("E2").Value = 0.04535142159914773
$ws.Range("F2").Value = 3.696417607249984
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 2.157564468342606
$ws.Range("J2").Value = 0.1037402433238022
$ws.Range("L2").Value = 0.3924357081747871
$ws.Range("M2").Value = 0.590738048185024
$ws.Range("N2").Value = 2.567041630754645
$ws.Range("B3").Value = 2.67863130748259
$ws.Range("C3").Value = 0.3305849875993658
$ws.Range("D3").Value = 0.008639454830692728
$ws.Range("E3").Value = 0.04498765152501871
$ws.Range("F3").Value = 3.688666605529122
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 2.160225284905664
$ws.Range("J3").Value = 0.1037164097846137
$ws.Range("L3").Value = 0.390246473843284
$ws.Range("M3").Value = 0.5774678069624812
$ws.Range("N3").Value = 2.590871131174303
$ws.Range("B4").Value = 2.625770047447702
$ws.Range("C4").Value = 0.3143289265103419
$ws.Range("D4").Value = 0.008297585275904851
$ws.Range("E4").Value = 0.04475934156084893
$ws.Range("F4").Value = 3.685739604632175
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 2.162872607125486
$ws.Range("J4").Value = 0.1037040203125503
$ws.Range("L4").Value = 0.3890673088526881
$ws.Range("M4").Value = 0.5696262734234168
$ws.Range("N4").Value = 2.606282021093939
$ws.Range("B5").Value = 2.604605926700913
$ws.Range("C5").Value = 0.3077545846118142
$ws.Range("D5").Value = 0.008157181559276694
$ws.Range("E5").Value = 0.04466504598544851
$ws.Range("F5").Value = 3.685006833174342
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 2.164205901958354
$ws.Range("J5").Value = 0.10369953747424
$ws.Range("L5").Value = 0.3886283279272718
$ws.Range("M5").Value = 0.5665078430750725
$ws.Range("N5").Value = 2.612757886682509
$ws.Range("B6").Value = 2.601114421896227
$ws.Range("C6").Value = 0.3066659361977031
$ws.Range("D6").Value = 0.008133800573130401
$ws.Range("E6").Value = 0.04464931198924749
$ws.Range("F6").Value = 3.684912919966536
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 2.164442654093065
$ws.Range("J6").Value = 0.1036988273336075
$ws.Range("L6").Value = 0.3885579454866104
$ws.Range("M6").Value = 0.5659946843121944
$ws.Range("N6").Value = 2.613845019354287
$ws.Range("B7").Value = 2.625483093406103
$ws.Range("C7").Value = 0.3142400601816746
$ws.Range("D7").Value = 0.008295696209629
$ws.Range("E7").Value = 0.04475807496399753
$ws.Range("F7").Value = 3.685727860625263
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 2.16288955850564
$ws.Range("J7").Value = 0.103703957561665
$ws.Range("L7").Value = 0.389061220361846
$ws.Range("M7").Value = 0.569583905175449
$ws.Range("N7").Value = 2.606368564372012
$ws.Range("B8").Value = 2.736329342626846
$ws.Range("C8").Value = 0.3481042622233019
$ws.Range("D8").Value = 0.00900066218813933
$ws.Range("E8").Value = 0.04522701464574386
$ws.Range("F8").Value = 3.693364361819292
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 2.15827131600917
$ws.Range("J8").Value = 0.1037315606173923
$ws.Range("L8").Value = 0.3916466290039224
$ws.Range("M8").Value = 0.5860988583692404
$ws.Range("N8").Value = 2.575096140284074
$ws.Range("B9").Value = 2.965633683675208
$ws.Range("C9").Value = 0.4161294350551543
$ws.Range("D9").Value = 0.01035213266977664
$ws.Range("E9").Value = 0.04610785708589837
$ws.Range("F9").Value = 3.72291425135495
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 2.157276658493046
$ws.Range("J9").Value = 0.1038034384504165
$ws.Range("L9").Value = 0.3980251022055654
$ws.Range("M9").Value = 0.6209188409757687
$ws.Range("N9").Value = 2.519964152420513
$ws.Range("B10").Value = 3.141436773030762
$ws.Range("C10").Value = 0.4671303913282259
$ws.Range("D10").Value = 0.01132903325892443
$ws.Range("E10").Value = 0.04673220311897985
$ws.Range("F10").Value = 3.753569981616067
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 2.161490814973746
$ws.Range("J10").Value = 0.1038670009281777
$ws.Range("L10").Value = 0.4035088698425682
$ws.Range("M10").Value = 0.647992586528197
$ws.Range("N10").Value = 2.483244440010118
$ws.Range("B11").Value = 3.223020637217132
$ws.Range("C11").Value = 0.4905624910068695
$ws.Range("D11").Value = 0.01177059343518749
$ws.Range("E11").Value = 0.04701145940365858
$ws.Range("F11").Value = 3.769472287082721
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 2.164488581522804
$ws.Range("J11").Value = 0.1038982398580206
$ws.Range("L11").Value = 0.406176810924336
$ws.Range("M11").Value = 0.6606349379672025
$ws.Range("N11").Value = 2.467364006691156
$ws.Range("B12").Value = 3.254146587534137
$ws.Range("C12").Value = 0.4994694771287413
$ws.Range("D12").Value = 0.01193744266496921
$ws.Range("E12").Value = 0.0471165352951024
$ws.Range("F12").Value = 3.775776478165767
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 2.16577968943912
$ws.Range("J12").Value = 0.1039104022066777
$ws.Range("L12").Value = 0.4072120051475707
$ws.Range("M12").Value = 0.6654692879338455
$ws.Range("N12").Value = 2.461469172150935
$ws.Range("B13").Value = 3.247432736640633
$ws.Range("C13").Value = 0.4975496893541731
$ws.Range("D13").Value = 0.01190152394343968
$ws.Range("E13").Value = 0.0470939350296451
$ws.Range("F13").Value = 3.77440618581241
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 2.165494683460807
$ws.Range("J13").Value = 0.103907768045314
$ws.Range("L13").Value = 0.4069879503574612
$ws.Range("M13").Value = 0.6644260347500079
$ws.Range("N13").Value = 2.462733444097523
$ws.Range("B14").Value = 3.225576737791414
$ws.Range("C14").Value = 0.4912945948377114
$ws.Range("D14").Value = 0.01178432716201883
$ws.Range("E14").Value = 0.04702011746539458
$ws.Range("F14").Value = 3.769985272855195
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 2.164591673938546
$ws.Range("J14").Value = 0.1038992337986269
$ws.Range("L14").Value = 0.4062614780188198
$ws.Range("M14").Value = 0.6610317213680403
$ws.Range("N14").Value = 2.466876653528288
$ws.Range("B15").Value = 3.212219520704991
$ws.Range("C15").Value = 0.4874675792045764
$ws.Range("D15").Value = 0.01171249529200935
$ws.Range("E15").Value = 0.04697481487721955
$ws.Range("F15").Value = 3.767314132318631
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 2.16405887511219
$ws.Range("J15").Value = 0.1038940496326468
$ws.Range("L15").Value = 0.4058197351275652
$ws.Range("M15").Value = 0.6589587238595271
$ws.Range("N15").Value = 2.469429961998628
$ws.Range("B16").Value = 3.136137507340266
$ws.Range("C16").Value = 0.4656037468787986
$ws.Range("D16").Value = 0.01130012326510865
$ws.Range("E16").Value = 0.04671385829661912
$ws.Range("F16").Value = 3.752570185415436
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 2.161316694366278
$ws.Range("J16").Value = 0.1038650060398734
$ws.Range("L16").Value = 0.4033379996943864
$ws.Range("M16").Value = 0.6471729493324432
$ws.Range("N16").Value = 2.484298862387092
$ws.Range("B17").Value = 3.089876312570254
$ws.Range("C17").Value = 0.4522506194153948
$ws.Range("D17").Value = 0.01104645375584212
$ws.Range("E17").Value = 0.04655255901755151
$ws.Range("F17").Value = 3.744027112404524
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 2.15991161609314
$ws.Range("J17").Value = 0.1038477830046691
$ws.Range("L17").Value = 0.4018599159921905
$ws.Range("M17").Value = 0.6400263721896238
$ws.Range("N17").Value = 2.49363160388566
$ws.Range("B18").Value = 3.063419628466761
$ws.Range("C18").Value = 0.4445920159601542
$ws.Range("D18").Value = 0.01090027986039388
$ws.Range("E18").Value = 0.0464593349858422
$ws.Range("F18").Value = 3.739297494425287
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 2.159205144083188
$ws.Range("J18").Value = 0.1038380955975615
$ws.Range("L18").Value = 0.4010260789926434
$ws.Range("M18").Value = 0.6359465670526916
$ws.Range("N18").Value = 2.499077046893525
$ws.Range("B19").Value = 3.054487864600389
$ws.Range("C19").Value = 0.4420026679857756
$ws.Range("D19").Value = 0.01085074008245002
$ws.Range("E19").Value = 0.04642769344794839
$ws.Range("F19").Value = 3.737727723508016
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 2.158983395022176
$ws.Range("J19").Value = 0.1038348532339537
$ws.Range("L19").Value = 0.4007465594537933
$ws.Range("M19").Value = 0.6345704898410034
$ws.Range("N19").Value = 2.500934080987889
$ws.Range("B20").Value = 3.094785210208443
$ws.Range("C20").Value = 0.4536698269420754
$ws.Range("D20").Value = 0.01107348489577831
$ws.Range("E20").Value = 0.04656977597612499
$ws.Range("F20").Value = 3.744917473505495
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 2.160050660718355
$ws.Range("J20").Value = 0.1038495937878388
$ws.Range("L20").Value = 0.4020155718946796
$ws.Range("M20").Value = 0.6407839582045085
$ws.Range("N20").Value = 2.492630094762262
$ws.Range("B21").Value = 3.231990077322621
$ws.Range("C21").Value = 0.4931309477616992
$ws.Range("D21").Value = 0.01181876008400806
$ws.Range("E21").Value = 0.04704181762333448
$ws.Range("F21").Value = 3.77127613213446
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 2.164852674005601
$ws.Range("J21").Value = 0.1039017314908275
$ws.Range("L21").Value = 0.4064741848838622
$ws.Range("M21").Value = 0.6620274385487832
$ws.Range("N21").Value = 2.465656467338214
$ws.Range("B22").Value = 3.323013341181309
$ws.Range("C22").Value = 0.5191180643464008
$ws.Range("D22").Value = 0.0123037645420716
$ws.Range("E22").Value = 0.04734641064539868
$ws.Range("F22").Value = 3.790149142400935
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 2.168900125759265
$ws.Range("J22").Value = 0.1039377459701338
$ws.Range("L22").Value = 0.4095332890089907
$ws.Range("M22").Value = 0.6761850534748319
$ws.Range("N22").Value = 2.448719907082328
$ws.Range("B23").Value = 3.274308677797706
$ws.Range("C23").Value = 0.5052300697462329
$ws.Range("D23").Value = 0.01204508254394909
$ws.Range("E23").Value = 0.04718419759767833
$ws.Range("F23").Value = 3.779925327765056
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 2.166656574021815
$ws.Range("J23").Value = 0.1039183473370713
$ws.Range("L23").Value = 0.4078873145604689
$ws.Range("M23").Value = 0.6686038050629648
$ws.Range("N23").Value = 2.457695819566474
$ws.Range("B24").Value = 3.092565463107292
$ws.Range("C24").Value = 0.4530281462970152
$ws.Range("D24").Value = 0.0110612651685571
$ws.Range("E24").Value = 0.04656199371960845
$ws.Range("F24").Value = 3.744514374686332
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 2.159987483052888
$ws.Range("J24").Value = 0.1038487744645877
$ws.Range("L24").Value = 0.4019451502310858
$ws.Range("M24").Value = 0.6404413636828963
$ws.Range("N24").Value = 2.493082628087826
$ws.Range("B25").Value = 2.902317230790175
$ws.Range("C25").Value = 0.3975498781724127
$ws.Range("D25").Value = 0.009989545790197951
$ws.Range("E25").Value = 0.04587363032347014
$ws.Range("F25").Value = 3.713353733863912
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 2.15667958081697
$ws.Range("J25").Value = 0.1037821004895676
$ws.Range("L25").Value = 0.3961594989409321
$ws.Range("M25").Value = 0.6112377519734977
$ws.Range("N25").Value = 2.534214764091153
